$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values scraped by the refresh job.
# D-column values that parse as plain numbers get a leading apostrophe so Excel
# keeps them as literal text (matching the sheets original inlineStr formatting)
# instead of silently converting them into floating point numbers.

$ws.Range("D2").Value = "54.407.19"
$ws.Range("E2").Value = "  -2.85%  "
$ws.Range("D3").Value = "2.287.40"
$ws.Range("E3").Value = "  -2.85%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'494.63"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("D6").Value = "'127.41"
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").Value = "2.286.34"
$ws.Range("E9").Value = "  -3.47%  "
$ws.Range("D10").Value = "'0.0943"
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("E13").Value = "  -3.64%  "
$ws.Range("D14").Value = "2.693.04"
$ws.Range("E14").Value = "  -2.92%  "
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "54.220.38"
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("D18").Value = "2.279.22"
$ws.Range("E18").Value = "  -4.60%  "
$ws.Range("D19").Value = "'9.96"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'4.06"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "'302.71"
$ws.Range("E21").Value = "  -2.68%  "
$ws.Range("D22").Value = "'6.42"
$ws.Range("E22").Value = "  +3.55%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "'5.38"
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("E25").Value = "  -2.62%  "
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("E28").Value = "  -2.84%  "
$ws.Range("E29").Value = "  +1.61%  "
$ws.Range("D30").Value = "'7.10"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "'168.09"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("E32").Value = "  -2.78%  "
$ws.Range("D33").Value = "0.0₃0684"
$ws.Range("E33").Value = "  -3.19%  "
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("D40").Value = "'0.872"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "'126.72"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").Value = "'238.90"
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("E51").Value = "  +0.27%  "
